$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Conversion" shifts right to E)
$ws.Columns("D").Insert()

# Try to match column C's width for the newly inserted column D
$ws.Columns.Item(4).ColumnWidth = 20.83

# Header for the new column
$ws.Range("D1").Value = "isMet"

# New column is all zeros for every data row
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
}

# Rename the AED variable used for row 17 (WQ_DIAG_TOT_TCHLA_SONDE)
$ws.Range("B17").Value = "WQ_DIAG_PHY_TCHLA"

# Match final selection from the saved workbook
$ws.Range("B17").Select()
